# Deploy edit: add a "2020" column (L) of data to the 9.c.1 indicator table
# and update the view/selection, matching the upstream OOXML diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column L values (2020 data) ---
$ws.Range("L4").Value = 2020
$ws.Range("L5").Value = 99
$ws.Range("L6").Value = 89.3
$ws.Range("L7").Value = 81.9

# --- Carry over the neighbouring cell formatting into the new column so the
#     new cells render the same way as the rest of the table (borders,
#     fonts, number formats, alignment). Source cells are picked so the
#     resulting style matches the one used by the rest of the row. ---
$ws.Range("K3").Copy() | Out-Null
$ws.Range("L3").PasteSpecial(-4122) | Out-Null     # xlPasteFormats (bottom border row)
$excel.CutCopyMode = $false

$ws.Range("K4").Copy() | Out-Null
$ws.Range("L4").PasteSpecial(-4122) | Out-Null     # xlPasteFormats (year header)
$excel.CutCopyMode = $false

$ws.Range("H5").Copy() | Out-Null
$ws.Range("L5").PasteSpecial(-4122) | Out-Null     # xlPasteFormats (0.0 numeric style)
$excel.CutCopyMode = $false

$ws.Range("K6").Copy() | Out-Null
$ws.Range("L6").PasteSpecial(-4122) | Out-Null     # xlPasteFormats (0.0 numeric style)
$excel.CutCopyMode = $false

$ws.Range("K7").Copy() | Out-Null
$ws.Range("L7").PasteSpecial(-4122) | Out-Null     # xlPasteFormats (0.0 numeric, bordered)
$excel.CutCopyMode = $false

# --- View: scroll so column B is the leftmost visible column, and leave the
#     active selection on N13 (matches the saved sheetView). ---
$ws.Activate()
$ws.Range("N13").Select()
$excel.ActiveWindow.ScrollColumn = 2
